# Generate Report for Handoff
#
# The "4345290d-7b98-49ac-89d3-937210843776" entry has now been handed off
# again (new "Latest Handoff Datetime"), so its status flips from
# "Handed back: in sync with en-US" to "Ready for handoff". The
# "b396911e-..." entry is now the most-recently handed-back entry, so it
# moves above the 4345290d entry (row 2) on every worksheet while the
# 4345290d entry drops to row 3 (Overview, zh-cn, de-de).
#
# Cells that carry hyperlinks must be rebuilt (delete + re-add) because
# this automation surface only supports appending new Hyperlink objects,
# not editing the TextToDisplay/Address of an existing one in place. The
# underlying hyperlink target URLs are unchanged by this edit, so each
# cell position keeps reusing the URL it already had.

$wb = $excel.ActiveWorkbook

function Set-PlainCell($ws, $addr, $text) {
    $ws.Range($addr).Value2 = $text
}

# ---------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Overview")

$ws1.Hyperlinks.Delete()
$ws1.Hyperlinks.Add($ws1.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/e2e/4345290d-7b98-49ac-89d3-937210843776.md", "", "", "b396911e-8a1e-4350-bc5f-2848b741994d.md")
$ws1.Hyperlinks.Add($ws1.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/e2e/b396911e-8a1e-4350-bc5f-2848b741994d.md", "", "", "4345290d-7b98-49ac-89d3-937210843776.md")
$ws1.Hyperlinks.Add($ws1.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/.localization-config", "", "", ".localization-config")

Set-PlainCell $ws1 "B2" "Handed back: in sync with en-US"
Set-PlainCell $ws1 "C2" "Handed back: in sync with en-US"

Set-PlainCell $ws1 "B3" "Ready for handoff"
Set-PlainCell $ws1 "C3" "Ready for handoff"

Set-PlainCell $ws1 "B4" "Not to be localized"
Set-PlainCell $ws1 "C4" "Not to be localized"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("zh-cn")

$ws2.Hyperlinks.Delete()
$ws2.Hyperlinks.Add($ws2.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/e2e/4345290d-7b98-49ac-89d3-937210843776.md", "", "", "b396911e-8a1e-4350-bc5f-2848b741994d.md")
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd021bb156a4f150e364a328c690786dd9542736/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.zh-cn.xlf", "", "", "b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6c8b1e487453c4b3f8af93daff336403bff719d6/e2e/4345290d-7b98-49ac-89d3-937210843776.md", "", "", "b396911e-8a1e-4350-bc5f-2848b741994d.md")
$ws2.Hyperlinks.Add($ws2.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7cddbd11b8bc1cc78e7da353e9ea0f759a1159d4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.zh-cn.xlf", "", "", "b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/e2e/b396911e-8a1e-4350-bc5f-2848b741994d.md", "", "", "4345290d-7b98-49ac-89d3-937210843776.md")
$ws2.Hyperlinks.Add($ws2.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/cd021bb156a4f150e364a328c690786dd9542736/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.zh-cn.xlf", "", "", "4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/6c8b1e487453c4b3f8af93daff336403bff719d6/e2e/b396911e-8a1e-4350-bc5f-2848b741994d.md", "", "", "4345290d-7b98-49ac-89d3-937210843776.md")
$ws2.Hyperlinks.Add($ws2.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/7cddbd11b8bc1cc78e7da353e9ea0f759a1159d4/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.zh-cn.xlf", "", "", "4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.zh-cn.xlf")
$ws2.Hyperlinks.Add($ws2.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/.localization-config", "", "", ".localization-config")

Set-PlainCell $ws2 "B2" "Handed back: in sync with en-US"
Set-PlainCell $ws2 "D2" "2016-03-08 08:29:05"
Set-PlainCell $ws2 "G2" "2016-03-08 08:29:33"
Set-PlainCell $ws2 "H2" "Include"

Set-PlainCell $ws2 "B3" "Ready for handoff"
Set-PlainCell $ws2 "D3" "2016-03-08 08:30:01"
Set-PlainCell $ws2 "G3" "2016-03-08 08:29:33"
Set-PlainCell $ws2 "H3" "Include"

Set-PlainCell $ws2 "B4" "Not to be localized"
Set-PlainCell $ws2 "D4" "0001-01-01 00:00:00"
Set-PlainCell $ws2 "G4" "0001-01-01 00:00:00"
Set-PlainCell $ws2 "H4" "Ignored"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("de-de")

$ws3.Hyperlinks.Delete()
$ws3.Hyperlinks.Add($ws3.Range("A2"), "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/e2e/4345290d-7b98-49ac-89d3-937210843776.md", "", "", "b396911e-8a1e-4350-bc5f-2848b741994d.md")
$ws3.Hyperlinks.Add($ws3.Range("C2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5201e81b08a588b71f0b1ced1d2f3b1d44edf0e4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.de-de.xlf", "", "", "b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E2"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/698f8a48cb1ab953b5b1f074db300f3ab9451e7e/e2e/4345290d-7b98-49ac-89d3-937210843776.md", "", "", "b396911e-8a1e-4350-bc5f-2848b741994d.md")
$ws3.Hyperlinks.Add($ws3.Range("F2"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/21f0d03a46a6ff0c1ca0091f68d609e0bcf45c05/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.de-de.xlf", "", "", "b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A3"), "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/e2e/b396911e-8a1e-4350-bc5f-2848b741994d.md", "", "", "4345290d-7b98-49ac-89d3-937210843776.md")
$ws3.Hyperlinks.Add($ws3.Range("C3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/5201e81b08a588b71f0b1ced1d2f3b1d44edf0e4/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.de-de.xlf", "", "", "4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("E3"), "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/698f8a48cb1ab953b5b1f074db300f3ab9451e7e/e2e/b396911e-8a1e-4350-bc5f-2848b741994d.md", "", "", "4345290d-7b98-49ac-89d3-937210843776.md")
$ws3.Hyperlinks.Add($ws3.Range("F3"), "https://github.com/OpenLocalizationTestOrg/olhandback/blob/21f0d03a46a6ff0c1ca0091f68d609e0bcf45c05/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/b396911e-8a1e-4350-bc5f-2848b741994d.d6cadb9cd9579b773674cb1118be25eac4eb70e1.de-de.xlf", "", "", "4345290d-7b98-49ac-89d3-937210843776.27c0bd8c516cd542ea3af7e11e0ce9c02e85b9ea.de-de.xlf")
$ws3.Hyperlinks.Add($ws3.Range("A4"), "https://github.com/OpenLocalizationTest/oltest/blob/147ea8c6b8413c786dcd49faf96b8423c3568fff/.localization-config", "", "", ".localization-config")

Set-PlainCell $ws3 "B2" "Handed back: in sync with en-US"
Set-PlainCell $ws3 "D2" "2016-03-08 08:29:11"
Set-PlainCell $ws3 "G2" "2016-03-08 08:29:40"
Set-PlainCell $ws3 "H2" "Include"

Set-PlainCell $ws3 "B3" "Ready for handoff"
Set-PlainCell $ws3 "D3" "2016-03-08 08:30:08"
Set-PlainCell $ws3 "G3" "2016-03-08 08:29:40"
Set-PlainCell $ws3 "H3" "Include"

Set-PlainCell $ws3 "B4" "Not to be localized"
Set-PlainCell $ws3 "D4" "0001-01-01 00:00:00"
Set-PlainCell $ws3 "G4" "0001-01-01 00:00:00"
Set-PlainCell $ws3 "H4" "Ignored"
